$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row above row 3 ("Zahlt Arbeitnehmer Rentenpauschale?") and
# shift the existing rows (old Arbeitgeberbeitrag.../U1/U2/.../Eintragsdatum
# block) down by one.
$ws.Rows.Item(3).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# New question row: mirrors the "kurzfristig beschaeftigt?" yes/no row above it.
$ws.Cells.Item(3, 1).Value = "Zahlt Arbeitnehmer Rentenpauschale?"
$ws.Cells.Item(3, 2).Value = "ja"

# The three "...beitrag..." labels become "...pauschale..." labels
# (values stay the same, now shifted to rows 4-6).
$ws.Cells.Item(4, 1).Value = "Arbeitgeberpauschale Krankenversicherung in Prozent"
$ws.Cells.Item(5, 1).Value = "Arbeitgeberpauschale Rentenversicherung in Prozent"
$ws.Cells.Item(6, 1).Value = "Arbeitnehmerpauschale Rentenversicherung in Prozent"

# Update the active selection to match the author's final cursor position.
$ws.Range("A8").Select()
